# Apply crypto price/volume updates per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "54.084.57"
$ws.Range("E2").Value = "  -10.88%  "

# Row 3
$ws.Range("D3").Value = "2.309.08"
$ws.Range("E3").Value = "  -20.51%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").Value = "  +0.12%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "445.00"
$ws.Range("E5").Value = "  -15.77%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "128.23"
$ws.Range("E6").Value = "  -11.25%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.998"
$ws.Range("E7").Value = "  -0.09%  "

# Row 8
$ws.Range("E8").Value = "  -14.69%  "

# Row 9
$ws.Range("D9").Value = "2.186.00"
$ws.Range("E9").Value = "  -24.96%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.37"
$ws.Range("E10").Value = "  -10.97%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0918"
$ws.Range("E11").Value = "  -15.52%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.310"
$ws.Range("E12").Value = "  -14.62%  "

# Row 13
$ws.Range("E13").Value = "  -3.50%  "

# Row 14
$ws.Range("D14").Value = "2.717.97"
$ws.Range("E14").Value = "  -20.32%  "

# Row 15
$ws.Range("D15").Value = "54.103.10"
$ws.Range("E15").Value = "  -10.80%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "18.78"
$ws.Range("E16").Value = "  -17.69%  "

# Row 17
$ws.Range("E17").Value = "  -14.50%  "

# Row 18
$ws.Range("D18").Value = "2.333.78"
$ws.Range("E18").Value = "  -19.73%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.05"
$ws.Range("E19").Value = "  -20.01%  "

# Row 20
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "298.28"
$ws.Range("E20").Value = "  -17.76%  "

# Row 21
$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.37"
$ws.Range("E21").Value = "  -20.16%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.998"
$ws.Range("E22").Value = "  -0.15%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.59"
$ws.Range("E23").Value = "  -1.65%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.33"
$ws.Range("E24").Value = "  -19.92%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "55.64"
$ws.Range("E25").Value = "  -14.18%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.999"
$ws.Range("E26").Value = "  +0.01%  "

# Row 27
$ws.Range("E27").Value = "  -13.68%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.368"
$ws.Range("E28").Value = "  -19.42%  "

# Row 29
$ws.Range("B29").Value = "InternetComputer(DFINITY)"
$ws.Range("C29").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.84"
$ws.Range("E29").Value = "  -13.42%  "

# Row 30
$ws.Range("B30").Value = "USDe"
$ws.Range("C30").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.995"
$ws.Range("E30").Value = "  -0.39%  "

# Row 31
$ws.Range("D31").Value = "0.0₃0708"
$ws.Range("E31").Value = "  -17.97%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "146.54"
$ws.Range("E32").Value = "  -3.59%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "16.98"
$ws.Range("E33").Value = "  -14.06%  "

# Row 34
$ws.Range("E34").Value = "  -19.92%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.68"
$ws.Range("E35").Value = "  -16.14%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.60"
$ws.Range("E36").Value = "  -18.22%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.836"
$ws.Range("E37").Value = "  -17.50%  "

# Row 38
$ws.Range("E38").Value = "  -16.72%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.995"
$ws.Range("E39").Value = "  -0.21%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "33.02"
$ws.Range("E40").Value = "  -12.29%  "

# Row 41
$ws.Range("E41").Value = "  -0.49%  "

# Row 42
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.14"
$ws.Range("E42").Value = "  -16.08%  "

# Row 43
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.22"
$ws.Range("E43").Value = "  -17.66%  "

# Row 44
$ws.Range("D44").Value = "1.927.85"
$ws.Range("E44").Value = "  -15.97%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0494"
$ws.Range("E45").Value = "  -15.29%  "

# Row 46
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.517"
$ws.Range("E46").Value = "  -20.38%  "

# Row 47
$ws.Range("B47").Value = "VeChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0208"
$ws.Range("E47").Value = "  -12.65%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0826"
$ws.Range("E48").Value = "  -10.87%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "16.08"
$ws.Range("E49").Value = "  -22.00%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "4.07"
$ws.Range("E50").Value = "  -18.99%  "

# Row 51
$ws.Range("E51").Value = "  -3.05%  "
